$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (37) with a new "inactivity" message entry, following the
# same pattern/formatting as the preceding rows (33-36).

# Copy formatting from the row above so the new cells pick up the same
# styles: A37 should look like A36 (style index 5), B37 should look like
# B36 (style index 6, the "Consolas" note style).
$ws.Range("A36").Copy()
$ws.Range("A37").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B36").Copy()
$ws.Range("B37").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new row's values.
$ws.Cells.Item(37, 1).Value = 6000
$ws.Cells.Item(37, 2).Value = "you are inactive"
$ws.Cells.Item(37, 4).Value = 4200
$ws.Cells.Item(37, 5).Value = 4200
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(37, 7).Value = "you are inactive low data"
$ws.Cells.Item(37, 9).Value = "message"

# Reflect the last-saved selection/cursor position in the sheet.
[void]$ws.Range("M15").Select()
